# Update mass-flow result values in both sheets with newly computed
# values from the updated input files (see commit message: "added new
# input files"). Only numeric result cells change; labels / structure
# stay the same.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Output_flows" ----
$ws1 = $wb.Worksheets.Item("Output_flows")

$ws1.Range("C7").Value = 0.0000001359819209298532
$ws1.Range("E7").Value = 0.00000008474995136565549
$ws1.Range("F7").Value = 0.00000002288248686872698
$ws1.Range("C12").Value = 0.00000000001818922729881697
$ws1.Range("E12").Value = 0.000000000005668165732666766
$ws1.Range("F12").Value = 0.000000000001530404747820027
$ws1.Range("C13").Value = 0.000000000004344442834704302
$ws1.Range("D13").Value = 0.0000000000003219712054650764
$ws1.Range("E13").Value = 0.000000000001353824524739553
$ws1.Range("F13").Value = 0.0000000000003655326216796793
$ws1.Range("C14").Value = 0.000000000002844024302994738
$ws1.Range("D14").Value = 0.0000000000008430944708420957
$ws1.Range("E14").Value = 0.0000000000008862609077492111
$ws1.Range("F14").Value = 0.000000000000239290445092287
$ws1.Range("C17").Value = 0.009542401001083108
$ws1.Range("E17").Value = 0.01189449325650512
$ws1.Range("F17").Value = 0.003211513179256382
$ws1.Range("C18").Value = 0.002310245342521586
$ws1.Range("D18").Value = 0.000006848672560848861
$ws1.Range("E18").Value = 0.002879694287043307
$ws1.Range("F18").Value = 0.0007775174575016929
$ws1.Range("C19").Value = 0.001854667993104582
$ws1.Range("D19").Value = 0.00002199222601888045
$ws1.Range("E19").Value = 0.002311822353151412
$ws1.Range("F19").Value = 0.0006241920353508812

# ---- Sheet "Input_flows" ----
$ws2 = $wb.Worksheets.Item("Input_flows")

$ws2.Range("C7").Value = 0.0000002436143591642356
$ws2.Range("C12").Value = 0.00000000002464427933841764
$ws2.Range("C13").Value = 0.000000000005964223951167562
$ws2.Range("C14").Value = 0.000000000004812670126678333
$ws2.Range("C17").Value = 0.02463056265127432
$ws2.Range("C18").Value = 0.005963309646617993
$ws2.Range("C19").Value = 0.004812674607625756

$wb.Save()
